# Updated symbol list on Sun Jan 29 10:41:57 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the crypto
# symbol rows. Values are written with a leading apostrophe so Excel keeps
# them as literal text (matching the original inline-string cells) instead
# of auto-converting them to numbers/percentages, and the style is reset
# back to "Normal" afterwards so no unintended number formatting/style is
# applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.03%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.12%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.149"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.95%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08191"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.89%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.978"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.47%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.158"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.75%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9273"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.07%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.80%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1967"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.38%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09038"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.46%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03508"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.05%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09817"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.08%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001393"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.09%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.006103"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.39%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.676"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.39%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.236"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'3.296"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-4.95%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.04%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1351"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'3.67%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.639"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-3.46%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.04377"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.20%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'2.03%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004801"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.71%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.33%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-10.19%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02162"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.08%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05213"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.28%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007410"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.93%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009803"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.36%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1374"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.09%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002125"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.79%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009874"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.02%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006387"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.20%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.28%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002764"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-9.90%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-37.63%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.28%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.28%"
$ws.Range("E51").Style = "Normal"
